$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B ("Loan Amount" -> shifts right to C, etc.)
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column B
$ws.Range("B1").Value = "ID"

# The inserted column copied formatting from column A (to its left); data rows
# B2:B4 should be unformatted like the other data columns, so clear that.
$ws.Range("B2:B4").ClearFormats()

# 2. Add headers for the three brand-new trailing columns (T, U, V), copying
#    the same bold/bordered header style used by the rest of row 1.
$ws.Range("S1").Copy()
$ws.Range("T1:V1").PasteSpecial(-4122)
$ws.Range("T1").Value = "Activity_Diversity"
$ws.Range("U1").Value = "Missed_Periods"
$ws.Range("V1").Value = "Cluster"

# 3. Write the full refreshed data block: columns B..V (ID + shifted metrics +
#    the 3 new metrics), rows 2..5 (row 5 is brand new -> cluster id 3).
$data = New-Object 'object[,]' 4,21

$data[0,0]=11322.875;      $data[0,1]=15034.73375;        $data[0,2]=19803.4025
$data[0,3]=0.2862888888888889; $data[0,4]=56.41666666666666; $data[0,5]=76854.97222222222
$data[0,6]=73894.44444444444;  $data[0,7]=3.266388888888889;  $data[0,8]=2.472222222222222
$data[0,9]=6.597222222222222;  $data[0,10]=4.736111111111111; $data[0,11]=4.639168055555555
$data[0,12]=0.1737996459739308; $data[0,13]=1.319563417002594; $data[0,14]=0.2163854196655239
$data[0,15]=0.02314480320969429; $data[0,16]=36.86111111111111; $data[0,17]=40.90277777777778
$data[0,18]=1.5; $data[0,19]=0.5763151622395932; $data[0,20]=0

$data[1,0]=11186.34426229508; $data[1,1]=15970.45180327869; $data[1,2]=20409.88508196721
$data[1,3]=0.2838344262295082; $data[1,4]=57.01639344262295; $data[1,5]=77559.40983606558
$data[1,6]=87053.27868852459; $data[1,7]=3.177540983606558; $data[1,8]=2.868852459016393
$data[1,9]=8.836065573770492; $data[1,10]=7.655737704918033; $data[1,11]=6.390165573770491
$data[1,12]=0.1618202756124067; $data[1,13]=1.285643004620153; $data[1,14]=0.2121999110056152
$data[1,15]=0.02622901379365879; $data[1,16]=46.9344262295082; $data[1,17]=41.29508196721311
$data[1,18]=1.80327868852459; $data[1,19]=0.6734740569407228; $data[1,20]=1

$data[2,0]=8555.448275862069; $data[2,1]=15814.80137931034; $data[2,2]=19380.13793103448
$data[2,3]=0.2774862068965517; $data[2,4]=54.44827586206897; $data[2,5]=95960.44827586207
$data[2,6]=72451.72413793103; $data[2,7]=3.114827586206896; $data[2,8]=1.931034482758621
$data[2,9]=8.793103448275861; $data[2,10]=19.06896551724138; $data[2,11]=5.428810344827586
$data[2,12]=0.1802418272197341; $data[2,13]=1.225321512882849; $data[2,14]=0.233650954215871
$data[2,15]=0.02294853230255155; $data[2,16]=59.20689655172414; $data[2,17]=35.41379310344828
$data[2,18]=1.586206896551724; $data[2,19]=0.5596866509069731; $data[2,20]=2

$data[3,0]=13751.18; $data[3,1]=12294.23; $data[3,2]=17074.1954
$data[3,3]=0.2851; $data[3,4]=56.54; $data[3,5]=64359.32
$data[3,6]=62152.62; $data[3,7]=3.2162; $data[3,8]=2.38
$data[3,9]=5.06; $data[3,10]=4.54; $data[3,11]=4.55717
$data[3,12]=0.1731970423347213; $data[3,13]=1.323189467343925; $data[3,14]=0.2166831423340825
$data[3,15]=0.02199783000512013; $data[3,16]=27.82; $data[3,17]=44.82
$data[3,18]=1.2; $data[3,19]=0.6121051984050072; $data[3,20]=3

$ws.Range("B2:V5").Value = $data

# 4. Column A row 5 (the new cluster id) needs the same bold/bordered style
#    used by the rest of column A.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3
